$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "grandes regiões" row (row 6) was just a bare label with no data -
# remove it entirely so the regional rows (norte/nordeste/sudeste/sul)
# shift up by one and pick up their correct data.
$ws.Rows("6:6").Delete()
